$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet is a "missing items" day-sale report. A new product row
# ("كريم فاتيكا 190 مل") is being inserted right where "محلول ملح" used to
# sit (row 45). The old "محلول ملح" row is pushed down to a brand-new row 46,
# the running total (now row 47) grows by the new product's selling price,
# the footer row moves from 47 -> 48, and the generated-at timestamp on that
# footer is refreshed.
# ---------------------------------------------------------------------------

# 1) Push everything from row 46 down one row, inserting a fresh blank row 46
#    (old row 46 "totals" -> 47, old row 47 "footer" -> 48).
$ws.Rows("46:46").Insert()

# 2) Give the new row 46 the same look as the data rows above it (row 45)
#    before filling in values.
$ws.Range("A45:Q45").Copy()
$ws.Range("A46:Q46").PasteSpecial(-4122)

# 3) Row 46 becomes the old "محلول ملح" line that used to live in row 45.
$ws.Range("A46").Value = 40
$ws.Range("B46").Value = ""
$ws.Range("C46").Value = "محلول ملح"
$ws.Range("H46").Value = "22:0"
$ws.Range("L46").Value = "'0"
$ws.Range("N46").Value = "24.00"
$ws.Range("P46").Value = "'24.0000"
$ws.Range("Q46").Value = "1:0"

# Re-merge the label cells for the newly inserted row, matching the layout
# used by every other data row.
$ws.Range("A46:B46").Merge()
$ws.Range("C46:G46").Merge()
$ws.Range("H46:K46").Merge()
$ws.Range("L46:M46").Merge()
$ws.Range("N46:O46").Merge()

# 4) Row 45 now holds the newly added product.
$ws.Range("C45").Value = "كريم فاتيكا 190 مل"
$ws.Range("H45").Value = "4:0"
$ws.Range("N45").Value = "65.00"
$ws.Range("P45").Value = "'65.0000"

# 5) The grand total (now on row 47) grows by the new item's selling price.
$ws.Range("P47").Value = 1990.7550000000001

# 6) The footer (now row 48) gets the refreshed "generated at" timestamp.
$ws.Range("A48").Value = "Friday, 1 August, 2025 6:31 PM"
